$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A4327').Value = 'Tử Vi tọa thủ tại cung đối Tử Tức'
$ws.Range('B4327').Value = 'Tử Vi tọa thủ tại cung đối Tử Tức'
$ws.Range('A4328').Value = 'Thiên Cơ tọa thủ tại cung đối Tử Tức'
$ws.Range('B4328').Value = 'Thiên Cơ tọa thủ tại cung đối Tử Tức'
$ws.Range('A4329').Value = 'Thái Dương tọa thủ tại cung đối Tử Tức'
$ws.Range('B4329').Value = 'Thái Dương tọa thủ tại cung đối Tử Tức'
$ws.Range('A4330').Value = 'Vũ Khúc tọa thủ tại cung đối Tử Tức'
$ws.Range('B4330').Value = 'Vũ Khúc tọa thủ tại cung đối Tử Tức'
$ws.Range('A4331').Value = 'Thiên Đồng tọa thủ tại cung đối Tử Tức'
$ws.Range('B4331').Value = 'Thiên Đồng tọa thủ tại cung đối Tử Tức'
$ws.Range('A4332').Value = 'Liêm Trinh tọa thủ tại cung đối Tử Tức'
$ws.Range('B4332').Value = 'Liêm Trinh tọa thủ tại cung đối Tử Tức'
$ws.Range('A4333').Value = 'Thiên Phủ tọa thủ tại cung đối Tử Tức'
$ws.Range('B4333').Value = 'Thiên Phủ tọa thủ tại cung đối Tử Tức'
$ws.Range('A4334').Value = 'Thái Âm tọa thủ tại cung đối Tử Tức'
$ws.Range('B4334').Value = 'Thái Âm tọa thủ tại cung đối Tử Tức'
$ws.Range('A4335').Value = 'Tham Lang tọa thủ tại cung đối Tử Tức'
$ws.Range('B4335').Value = 'Tham Lang tọa thủ tại cung đối Tử Tức'
$ws.Range('A4336').Value = 'Cự Môn tọa thủ tại cung đối Tử Tức'
$ws.Range('B4336').Value = 'Cự Môn tọa thủ tại cung đối Tử Tức'
$ws.Range('A4337').Value = 'Thiên Tướng tọa thủ tại cung đối Tử Tức'
$ws.Range('B4337').Value = 'Thiên Tướng tọa thủ tại cung đối Tử Tức'
$ws.Range('A4338').Value = 'Thiên Lương tọa thủ tại cung đối Tử Tức'
$ws.Range('B4338').Value = 'Thiên Lương tọa thủ tại cung đối Tử Tức'
$ws.Range('A4339').Value = 'Thất Sát tọa thủ tại cung đối Tử Tức'
$ws.Range('B4339').Value = 'Thất Sát tọa thủ tại cung đối Tử Tức'
$ws.Range('A4340').Value = 'Phá Quân tọa thủ tại cung đối Tử Tức'
$ws.Range('B4340').Value = 'Phá Quân tọa thủ tại cung đối Tử Tức'
$ws.Range('A4342').Value = 'Tử Vi đồng cung Thiên Cơ tại cung đối Tử Tức'
$ws.Range('B4342').Value = 'Tử Vi đồng cung Thiên Cơ tại cung đối Tử Tức'
$ws.Range('A4343').Value = 'Tử Vi đồng cung Thái Dương tại cung đối Tử Tức'
$ws.Range('B4343').Value = 'Tử Vi đồng cung Thái Dương tại cung đối Tử Tức'
$ws.Range('A4344').Value = 'Tử Vi đồng cung Vũ Khúc tại cung đối Tử Tức'
$ws.Range('B4344').Value = 'Tử Vi đồng cung Vũ Khúc tại cung đối Tử Tức'
$ws.Range('A4345').Value = 'Tử Vi đồng cung Thiên Đồng tại cung đối Tử Tức'
$ws.Range('B4345').Value = 'Tử Vi đồng cung Thiên Đồng tại cung đối Tử Tức'
$ws.Range('A4346').Value = 'Tử Vi đồng cung Liêm Trinh tại cung đối Tử Tức'
$ws.Range('B4346').Value = 'Tử Vi đồng cung Liêm Trinh tại cung đối Tử Tức'
$ws.Range('A4347').Value = 'Tử Vi đồng cung Thiên Phủ tại cung đối Tử Tức'
$ws.Range('B4347').Value = 'Tử Vi đồng cung Thiên Phủ tại cung đối Tử Tức'
$ws.Range('A4348').Value = 'Tử Vi đồng cung Thái Âm tại cung đối Tử Tức'
$ws.Range('B4348').Value = 'Tử Vi đồng cung Thái Âm tại cung đối Tử Tức'
$ws.Range('A4349').Value = 'Tử Vi đồng cung Tham Lang tại cung đối Tử Tức'
$ws.Range('B4349').Value = 'Tử Vi đồng cung Tham Lang tại cung đối Tử Tức'
$ws.Range('A4350').Value = 'Tử Vi đồng cung Cự Môn tại cung đối Tử Tức'
$ws.Range('B4350').Value = 'Tử Vi đồng cung Cự Môn tại cung đối Tử Tức'
$ws.Range('A4351').Value = 'Tử Vi đồng cung Thiên Tướng tại cung đối Tử Tức'
$ws.Range('B4351').Value = 'Tử Vi đồng cung Thiên Tướng tại cung đối Tử Tức'
$ws.Range('A4352').Value = 'Tử Vi đồng cung Thiên Lương tại cung đối Tử Tức'
$ws.Range('B4352').Value = 'Tử Vi đồng cung Thiên Lương tại cung đối Tử Tức'
$ws.Range('A4353').Value = 'Tử Vi đồng cung Thất Sát tại cung đối Tử Tức'
$ws.Range('B4353').Value = 'Tử Vi đồng cung Thất Sát tại cung đối Tử Tức'
$ws.Range('A4354').Value = 'Tử Vi đồng cung Phá Quân tại cung đối Tử Tức'
$ws.Range('B4354').Value = 'Tử Vi đồng cung Phá Quân tại cung đối Tử Tức'
$ws.Range('A4355').Value = 'Thiên Cơ đồng cung Thái Dương tại cung đối Tử Tức'
$ws.Range('B4355').Value = 'Thiên Cơ đồng cung Thái Dương tại cung đối Tử Tức'
$ws.Range('A4356').Value = 'Thiên Cơ đồng cung Vũ Khúc tại cung đối Tử Tức'
$ws.Range('B4356').Value = 'Thiên Cơ đồng cung Vũ Khúc tại cung đối Tử Tức'
$ws.Range('A4357').Value = 'Thiên Cơ đồng cung Thiên Đồng tại cung đối Tử Tức'
$ws.Range('B4357').Value = 'Thiên Cơ đồng cung Thiên Đồng tại cung đối Tử Tức'
$ws.Range('A4358').Value = 'Thiên Cơ đồng cung Liêm Trinh tại cung đối Tử Tức'
$ws.Range('B4358').Value = 'Thiên Cơ đồng cung Liêm Trinh tại cung đối Tử Tức'
$ws.Range('A4359').Value = 'Thiên Cơ đồng cung Thiên Phủ tại cung đối Tử Tức'
$ws.Range('B4359').Value = 'Thiên Cơ đồng cung Thiên Phủ tại cung đối Tử Tức'
$ws.Range('A4360').Value = 'Thiên Cơ đồng cung Thái Âm tại cung đối Tử Tức'
$ws.Range('B4360').Value = 'Thiên Cơ đồng cung Thái Âm tại cung đối Tử Tức'
$ws.Range('A4361').Value = 'Thiên Cơ đồng cung Tham Lang tại cung đối Tử Tức'
$ws.Range('B4361').Value = 'Thiên Cơ đồng cung Tham Lang tại cung đối Tử Tức'
$ws.Range('A4362').Value = 'Thiên Cơ đồng cung Cự Môn tại cung đối Tử Tức'
$ws.Range('B4362').Value = 'Thiên Cơ đồng cung Cự Môn tại cung đối Tử Tức'
$ws.Range('A4363').Value = 'Thiên Cơ đồng cung Thiên Tướng tại cung đối Tử Tức'
$ws.Range('B4363').Value = 'Thiên Cơ đồng cung Thiên Tướng tại cung đối Tử Tức'
$ws.Range('A4364').Value = 'Thiên Cơ đồng cung Thiên Lương tại cung đối Tử Tức'
$ws.Range('B4364').Value = 'Thiên Cơ đồng cung Thiên Lương tại cung đối Tử Tức'
$ws.Range('A4365').Value = 'Thiên Cơ đồng cung Thất Sát tại cung đối Tử Tức'
$ws.Range('B4365').Value = 'Thiên Cơ đồng cung Thất Sát tại cung đối Tử Tức'
$ws.Range('A4366').Value = 'Thiên Cơ đồng cung Phá Quân tại cung đối Tử Tức'
$ws.Range('B4366').Value = 'Thiên Cơ đồng cung Phá Quân tại cung đối Tử Tức'
$ws.Range('A4367').Value = 'Thái Dương đồng cung Vũ Khúc tại cung đối Tử Tức'
$ws.Range('B4367').Value = 'Thái Dương đồng cung Vũ Khúc tại cung đối Tử Tức'
$ws.Range('A4368').Value = 'Thái Dương đồng cung Thiên Đồng tại cung đối Tử Tức'
$ws.Range('B4368').Value = 'Thái Dương đồng cung Thiên Đồng tại cung đối Tử Tức'
$ws.Range('A4369').Value = 'Thái Dương đồng cung Liêm Trinh tại cung đối Tử Tức'
$ws.Range('B4369').Value = 'Thái Dương đồng cung Liêm Trinh tại cung đối Tử Tức'
$ws.Range('A4370').Value = 'Thái Dương đồng cung Thiên Phủ tại cung đối Tử Tức'
$ws.Range('B4370').Value = 'Thái Dương đồng cung Thiên Phủ tại cung đối Tử Tức'
$ws.Range('A4371').Value = 'Thái Dương đồng cung Thái Âm tại cung đối Tử Tức'
$ws.Range('B4371').Value = 'Thái Dương đồng cung Thái Âm tại cung đối Tử Tức'
$ws.Range('A4372').Value = 'Thái Dương đồng cung Tham Lang tại cung đối Tử Tức'
$ws.Range('B4372').Value = 'Thái Dương đồng cung Tham Lang tại cung đối Tử Tức'
$ws.Range('A4373').Value = 'Thái Dương đồng cung Cự Môn tại cung đối Tử Tức'
$ws.Range('B4373').Value = 'Thái Dương đồng cung Cự Môn tại cung đối Tử Tức'
$ws.Range('A4374').Value = 'Thái Dương đồng cung Thiên Tướng tại cung đối Tử Tức'
$ws.Range('B4374').Value = 'Thái Dương đồng cung Thiên Tướng tại cung đối Tử Tức'
$ws.Range('A4375').Value = 'Thái Dương đồng cung Thiên Lương tại cung đối Tử Tức'
$ws.Range('B4375').Value = 'Thái Dương đồng cung Thiên Lương tại cung đối Tử Tức'
$ws.Range('A4376').Value = 'Thái Dương đồng cung Thất Sát tại cung đối Tử Tức'
$ws.Range('B4376').Value = 'Thái Dương đồng cung Thất Sát tại cung đối Tử Tức'
$ws.Range('A4377').Value = 'Thái Dương đồng cung Phá Quân tại cung đối Tử Tức'
$ws.Range('B4377').Value = 'Thái Dương đồng cung Phá Quân tại cung đối Tử Tức'
$ws.Range('A4378').Value = 'Vũ Khúc đồng cung Thiên Đồng tại cung đối Tử Tức'
$ws.Range('B4378').Value = 'Vũ Khúc đồng cung Thiên Đồng tại cung đối Tử Tức'
$ws.Range('A4379').Value = 'Vũ Khúc đồng cung Liêm Trinh tại cung đối Tử Tức'
$ws.Range('B4379').Value = 'Vũ Khúc đồng cung Liêm Trinh tại cung đối Tử Tức'
$ws.Range('A4380').Value = 'Vũ Khúc đồng cung Thiên Phủ tại cung đối Tử Tức'
$ws.Range('B4380').Value = 'Vũ Khúc đồng cung Thiên Phủ tại cung đối Tử Tức'
$ws.Range('A4381').Value = 'Vũ Khúc đồng cung Thái Âm tại cung đối Tử Tức'
$ws.Range('B4381').Value = 'Vũ Khúc đồng cung Thái Âm tại cung đối Tử Tức'
$ws.Range('A4382').Value = 'Vũ Khúc đồng cung Tham Lang tại cung đối Tử Tức'
$ws.Range('B4382').Value = 'Vũ Khúc đồng cung Tham Lang tại cung đối Tử Tức'
$ws.Range('A4383').Value = 'Vũ Khúc đồng cung Cự Môn tại cung đối Tử Tức'
$ws.Range('B4383').Value = 'Vũ Khúc đồng cung Cự Môn tại cung đối Tử Tức'
$ws.Range('A4384').Value = 'Vũ Khúc đồng cung Thiên Tướng tại cung đối Tử Tức'
$ws.Range('B4384').Value = 'Vũ Khúc đồng cung Thiên Tướng tại cung đối Tử Tức'
$ws.Range('A4385').Value = 'Vũ Khúc đồng cung Thiên Lương tại cung đối Tử Tức'
$ws.Range('B4385').Value = 'Vũ Khúc đồng cung Thiên Lương tại cung đối Tử Tức'
$ws.Range('A4386').Value = 'Vũ Khúc đồng cung Thất Sát tại cung đối Tử Tức'
$ws.Range('B4386').Value = 'Vũ Khúc đồng cung Thất Sát tại cung đối Tử Tức'
$ws.Range('A4387').Value = 'Vũ Khúc đồng cung Phá Quân tại cung đối Tử Tức'
$ws.Range('B4387').Value = 'Vũ Khúc đồng cung Phá Quân tại cung đối Tử Tức'
$ws.Range('A4388').Value = 'Thiên Đồng đồng cung Liêm Trinh tại cung đối Tử Tức'
$ws.Range('B4388').Value = 'Thiên Đồng đồng cung Liêm Trinh tại cung đối Tử Tức'
$ws.Range('A4389').Value = 'Thiên Đồng đồng cung Thiên Phủ tại cung đối Tử Tức'
$ws.Range('B4389').Value = 'Thiên Đồng đồng cung Thiên Phủ tại cung đối Tử Tức'
$ws.Range('A4390').Value = 'Thiên Đồng đồng cung Thái Âm tại cung đối Tử Tức'
$ws.Range('B4390').Value = 'Thiên Đồng đồng cung Thái Âm tại cung đối Tử Tức'
$ws.Range('A4391').Value = 'Thiên Đồng đồng cung Tham Lang tại cung đối Tử Tức'
$ws.Range('B4391').Value = 'Thiên Đồng đồng cung Tham Lang tại cung đối Tử Tức'
$ws.Range('A4392').Value = 'Thiên Đồng đồng cung Cự Môn tại cung đối Tử Tức'
$ws.Range('B4392').Value = 'Thiên Đồng đồng cung Cự Môn tại cung đối Tử Tức'
$ws.Range('A4393').Value = 'Thiên Đồng đồng cung Thiên Tướng tại cung đối Tử Tức'
$ws.Range('B4393').Value = 'Thiên Đồng đồng cung Thiên Tướng tại cung đối Tử Tức'
$ws.Range('A4394').Value = 'Thiên Đồng đồng cung Thiên Lương tại cung đối Tử Tức'
$ws.Range('B4394').Value = 'Thiên Đồng đồng cung Thiên Lương tại cung đối Tử Tức'
$ws.Range('A4395').Value = 'Thiên Đồng đồng cung Thất Sát tại cung đối Tử Tức'
$ws.Range('B4395').Value = 'Thiên Đồng đồng cung Thất Sát tại cung đối Tử Tức'
$ws.Range('A4396').Value = 'Thiên Đồng đồng cung Phá Quân tại cung đối Tử Tức'
$ws.Range('B4396').Value = 'Thiên Đồng đồng cung Phá Quân tại cung đối Tử Tức'
$ws.Range('A4397').Value = 'Liêm Trinh đồng cung Thiên Phủ tại cung đối Tử Tức'
$ws.Range('B4397').Value = 'Liêm Trinh đồng cung Thiên Phủ tại cung đối Tử Tức'
$ws.Range('A4398').Value = 'Liêm Trinh đồng cung Thái Âm tại cung đối Tử Tức'
$ws.Range('B4398').Value = 'Liêm Trinh đồng cung Thái Âm tại cung đối Tử Tức'
$ws.Range('A4399').Value = 'Liêm Trinh đồng cung Tham Lang tại cung đối Tử Tức'
$ws.Range('B4399').Value = 'Liêm Trinh đồng cung Tham Lang tại cung đối Tử Tức'
$ws.Range('A4400').Value = 'Liêm Trinh đồng cung Cự Môn tại cung đối Tử Tức'
$ws.Range('B4400').Value = 'Liêm Trinh đồng cung Cự Môn tại cung đối Tử Tức'
$ws.Range('A4401').Value = 'Liêm Trinh đồng cung Thiên Tướng tại cung đối Tử Tức'
$ws.Range('B4401').Value = 'Liêm Trinh đồng cung Thiên Tướng tại cung đối Tử Tức'
$ws.Range('A4402').Value = 'Liêm Trinh đồng cung Thiên Lương tại cung đối Tử Tức'
$ws.Range('B4402').Value = 'Liêm Trinh đồng cung Thiên Lương tại cung đối Tử Tức'
$ws.Range('A4403').Value = 'Liêm Trinh đồng cung Thất Sát tại cung đối Tử Tức'
$ws.Range('B4403').Value = 'Liêm Trinh đồng cung Thất Sát tại cung đối Tử Tức'
$ws.Range('A4404').Value = 'Liêm Trinh đồng cung Phá Quân tại cung đối Tử Tức'
$ws.Range('B4404').Value = 'Liêm Trinh đồng cung Phá Quân tại cung đối Tử Tức'
$ws.Range('A4405').Value = 'Thiên Phủ đồng cung Thái Âm tại cung đối Tử Tức'
$ws.Range('B4405').Value = 'Thiên Phủ đồng cung Thái Âm tại cung đối Tử Tức'
$ws.Range('A4406').Value = 'Thiên Phủ đồng cung Tham Lang tại cung đối Tử Tức'
$ws.Range('B4406').Value = 'Thiên Phủ đồng cung Tham Lang tại cung đối Tử Tức'
$ws.Range('A4407').Value = 'Thiên Phủ đồng cung Cự Môn tại cung đối Tử Tức'
$ws.Range('B4407').Value = 'Thiên Phủ đồng cung Cự Môn tại cung đối Tử Tức'
$ws.Range('A4408').Value = 'Thiên Phủ đồng cung Thiên Tướng tại cung đối Tử Tức'
$ws.Range('B4408').Value = 'Thiên Phủ đồng cung Thiên Tướng tại cung đối Tử Tức'
$ws.Range('A4409').Value = 'Thiên Phủ đồng cung Thiên Lương tại cung đối Tử Tức'
$ws.Range('B4409').Value = 'Thiên Phủ đồng cung Thiên Lương tại cung đối Tử Tức'
$ws.Range('A4410').Value = 'Thiên Phủ đồng cung Thất Sát tại cung đối Tử Tức'
$ws.Range('B4410').Value = 'Thiên Phủ đồng cung Thất Sát tại cung đối Tử Tức'
$ws.Range('A4411').Value = 'Thiên Phủ đồng cung Phá Quân tại cung đối Tử Tức'
$ws.Range('B4411').Value = 'Thiên Phủ đồng cung Phá Quân tại cung đối Tử Tức'
$ws.Range('A4412').Value = 'Thái Âm đồng cung Tham Lang tại cung đối Tử Tức'
$ws.Range('B4412').Value = 'Thái Âm đồng cung Tham Lang tại cung đối Tử Tức'
$ws.Range('A4413').Value = 'Thái Âm đồng cung Cự Môn tại cung đối Tử Tức'
$ws.Range('B4413').Value = 'Thái Âm đồng cung Cự Môn tại cung đối Tử Tức'
$ws.Range('A4414').Value = 'Thái Âm đồng cung Thiên Tướng tại cung đối Tử Tức'
$ws.Range('B4414').Value = 'Thái Âm đồng cung Thiên Tướng tại cung đối Tử Tức'
$ws.Range('A4415').Value = 'Thái Âm đồng cung Thiên Lương tại cung đối Tử Tức'
$ws.Range('B4415').Value = 'Thái Âm đồng cung Thiên Lương tại cung đối Tử Tức'
$ws.Range('A4416').Value = 'Thái Âm đồng cung Thất Sát tại cung đối Tử Tức'
$ws.Range('B4416').Value = 'Thái Âm đồng cung Thất Sát tại cung đối Tử Tức'
$ws.Range('A4417').Value = 'Thái Âm đồng cung Phá Quân tại cung đối Tử Tức'
$ws.Range('B4417').Value = 'Thái Âm đồng cung Phá Quân tại cung đối Tử Tức'
$ws.Range('A4418').Value = 'Tham Lang đồng cung Cự Môn tại cung đối Tử Tức'
$ws.Range('B4418').Value = 'Tham Lang đồng cung Cự Môn tại cung đối Tử Tức'
$ws.Range('A4419').Value = 'Tham Lang đồng cung Thiên Tướng tại cung đối Tử Tức'
$ws.Range('B4419').Value = 'Tham Lang đồng cung Thiên Tướng tại cung đối Tử Tức'
$ws.Range('A4420').Value = 'Tham Lang đồng cung Thiên Lương tại cung đối Tử Tức'
$ws.Range('B4420').Value = 'Tham Lang đồng cung Thiên Lương tại cung đối Tử Tức'
$ws.Range('A4421').Value = 'Tham Lang đồng cung Thất Sát tại cung đối Tử Tức'
$ws.Range('B4421').Value = 'Tham Lang đồng cung Thất Sát tại cung đối Tử Tức'
$ws.Range('A4422').Value = 'Tham Lang đồng cung Phá Quân tại cung đối Tử Tức'
$ws.Range('B4422').Value = 'Tham Lang đồng cung Phá Quân tại cung đối Tử Tức'
$ws.Range('A4423').Value = 'Cự Môn đồng cung Thiên Tướng tại cung đối Tử Tức'
$ws.Range('B4423').Value = 'Cự Môn đồng cung Thiên Tướng tại cung đối Tử Tức'
$ws.Range('A4424').Value = 'Cự Môn đồng cung Thiên Lương tại cung đối Tử Tức'
$ws.Range('B4424').Value = 'Cự Môn đồng cung Thiên Lương tại cung đối Tử Tức'
$ws.Range('A4425').Value = 'Cự Môn đồng cung Thất Sát tại cung đối Tử Tức'
$ws.Range('B4425').Value = 'Cự Môn đồng cung Thất Sát tại cung đối Tử Tức'
$ws.Range('A4426').Value = 'Cự Môn đồng cung Phá Quân tại cung đối Tử Tức'
$ws.Range('B4426').Value = 'Cự Môn đồng cung Phá Quân tại cung đối Tử Tức'
$ws.Range('A4427').Value = 'Thiên Tướng đồng cung Thiên Lương tại cung đối Tử Tức'
$ws.Range('B4427').Value = 'Thiên Tướng đồng cung Thiên Lương tại cung đối Tử Tức'
$ws.Range('A4428').Value = 'Thiên Tướng đồng cung Thất Sát tại cung đối Tử Tức'
$ws.Range('B4428').Value = 'Thiên Tướng đồng cung Thất Sát tại cung đối Tử Tức'
$ws.Range('A4429').Value = 'Thiên Tướng đồng cung Phá Quân tại cung đối Tử Tức'
$ws.Range('B4429').Value = 'Thiên Tướng đồng cung Phá Quân tại cung đối Tử Tức'
$ws.Range('A4430').Value = 'Thiên Lương đồng cung Thất Sát tại cung đối Tử Tức'
$ws.Range('B4430').Value = 'Thiên Lương đồng cung Thất Sát tại cung đối Tử Tức'
$ws.Range('A4431').Value = 'Thiên Lương đồng cung Phá Quân tại cung đối Tử Tức'
$ws.Range('B4431').Value = 'Thiên Lương đồng cung Phá Quân tại cung đối Tử Tức'
$ws.Range('A4432').Value = 'Thất Sát đồng cung Phá Quân tại cung đối Tử Tức'
$ws.Range('B4432').Value = 'Thất Sát đồng cung Phá Quân tại cung đối Tử Tức'

# Update selection to match the final saved state (B4342:B4432)
$null = $ws.Range('B4342:B4432').Select()
